# Updated link for slides.
# The worksheet physically holding the Date/Topic/.../Slides/Homework/Lab
# schedule data is named "Meetups" in this workbook (tab-name/data is
# swapped vs. the other sheet, "Schedule").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meetups")
$ws.Activate()

# Row 7 ("Inference for Numerical Data"): the slides link pointed at the
# (wrong) 06- prefixed file; it should point at the 07- prefixed one.
$ws.Range("F7").Value = "/slides/07-Inference_for_Numerical_Data.html"

# Row 16 ("Presentations") had the bayesian chapter link in its Prepare
# column by mistake; it belongs on row 14 ("Bayesian Analysis").
$ws.Range("E14").Value = $ws.Range("E16").Value2
$ws.Range("E16").ClearContents()

# Update the active cell/selection to match.
$ws.Range("F5").Select() | Out-Null
